$d = $word.ActiveDocument

# Recolor the sticky-note title rectangles that are filled with #c1106e
# to the new color #36a98b (keeps every other shape untouched).
#
# Quirk of this COM-interop runtime: Shapes.Item(i) writes shape
# formatting properties using *document (XML) order* indexing, which can
# differ from the order its getters (Name, etc.) appear to enumerate in.
# So rather than trust Shapes.Item(i).Name while looping 1..Count, we
# first recover the true document order (and each shape's current fill
# color) from Content.WordOpenXML, then issue the Fill.ForeColor.RGB
# writes using that document-order index - this lands reliably on the
# intended shape.

$oldColor = 0x6E10C1  # BGR-long for #c1106e
$newColor = 0x8BA936  # BGR-long for #36a98b

$xml = $d.Content.WordOpenXML
$shapeMatches = [regex]::Matches($xml, 'docPr id="(\d+)" name="([^"]+)"')

$targetIndexes = @()
$pos = 0
foreach ($m in $shapeMatches) {
    $pos = $pos + 1
    $startPos = $m.Index
    $snippetLen = [Math]::Min(800, $xml.Length - $startPos)
    $snippet = $xml.Substring($startPos, $snippetLen)
    $colorMatch = [regex]::Match($snippet, 'srgbClr val="([0-9a-fA-F]{6})"')
    if ($colorMatch.Success -and ($colorMatch.Groups[1].Value -ieq "c1106e")) {
        $targetIndexes += $pos
    }
}

$changed = 0
foreach ($idx in $targetIndexes) {
    $s = $d.Shapes.Item($idx)
    $s.Fill.ForeColor.RGB = $newColor
    $changed = $changed + 1
}

Write-Output "Shapes recolored: $changed (indexes: $($targetIndexes -join ','))"
